# Append: 2025-11-25 12:38 JST
# Two new scraped job listings are merged into the existing list (rows 4
# and 9), pushing the later rows down. Every row's "取得日時" (fetched-at)
# timestamp is refreshed to the new run's time. We rewrite the whole data
# block (rows 2-12) explicitly with the final values, which is equivalent
# to the row-insert + shift the original diff shows but is unambiguous
# here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ts = "2025-11-25 12:38:07"

# --- Row 2 (unchanged job, only the timestamp refreshes) -------------------
$ws.Range("A2").Value = $ts

# --- Row 3 (unchanged job, only the timestamp refreshes) -------------------
$ws.Range("A3").Value = $ts

# --- Row 4: NEW job inserted here --------------------------------------
$ws.Range("A4").Value = $ts
$ws.Range("B4").Value = "【急募】古いPHPとPerlプログラムのアップデート依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5440861"
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = "○PHP"

# --- Row 5: was row 4 (シティヘブン...) --------------------------------
$ws.Range("A5").Value = $ts
$ws.Range("B5").Value = "【急募】シティヘブンの出勤情報を自動取得・管理したい!"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5440436"
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = "◇管理"

# --- Row 6: was row 5 (進行管理...) ------------------------------------
$ws.Range("A6").Value = $ts
$ws.Range("B6").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "~ 5,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = "◇管理"

# --- Row 7: was row 6 (初回 n8n+Gemini...) -----------------------------
$ws.Range("A7").Value = $ts
$ws.Range("B7").Value = "初回 n8n+Gemini+Typefully+GoogleスプレッドのX/Threads自動投稿システム"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5440440"
$ws.Range("G7").Value = 33
$ws.Range("H7").ClearContents()

# --- Row 8: was row 7 (急募 限定公開...) -------------------------------
$ws.Range("A8").Value = $ts
$ws.Range("B8").Value = "急募 限定公開 PR 限定公開の仕事"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5440230"
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

# --- Row 9: NEW job inserted here ---------------------------------------
$ws.Range("A9").Value = $ts
$ws.Range("B9").Value = "〖リモート可〗Delphiエンジニア募集"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5341051"
$ws.Range("G9").Value = 25
$ws.Range("H9").ClearContents()

# --- Row 10: was row 8 (リーダー募集...) -------------------------------
$ws.Range("A10").Value = $ts
$ws.Range("B10").Value = "【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5435080"
$ws.Range("G10").Value = 25
$ws.Range("H10").ClearContents()

# --- Row 11: was row 9 (若手歓迎...) ------------------------------------
$ws.Range("A11").Value = $ts
$ws.Range("B11").Value = "【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5435079"
$ws.Range("G11").Value = 25
$ws.Range("H11").ClearContents()

# --- Row 12: was row 10 (保守運用...) -----------------------------------
$ws.Range("A12").Value = $ts
$ws.Range("B12").Value = "【急募】弊社Websiteの保守運用をお任せできる方を探しています!"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5440806"
$ws.Range("G12").Value = 18
$ws.Range("H12").ClearContents()

# --- Rebuild the hyperlinks on column F (rows 2-12) ----------------------
# Deleting hyperlinks from any single cell in the sheet clears the whole
# worksheet's Hyperlinks collection in this host, so do that once, then
# re-add a fresh hyperlink per row in order.
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5440461")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5440861")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5440436")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5440440")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5440230")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5341051")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5435080")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5435079")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5440806")
